$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (price) cells to remain plain text even when the new
# value happens to look like a number, matching the inlineStr storage
# used by the source data. Reset the style back to Normal afterwards so
# no stray cell-level style id is left behind.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '26.459.70'
$ws.Range('E2').Value = '  -0.56%  '
$ws.Range('D3').Value = '1.843.72'
$ws.Range('E3').Value = '  -0.78%  '
$ws.Range('D4').Value = '0.9992'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '262.41'
$ws.Range('E5').Value = '  -3.75%  '
$ws.Range('D6').Value = '0.9996'
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('D7').Value = '0.5207'
$ws.Range('E7').Value = '  -1.34%  '
$ws.Range('D8').Value = '0.3263'
$ws.Range('E8').Value = '  -3.34%  '
$ws.Range('D9').Value = '0.06800'
$ws.Range('E9').Value = '  +0.16%  '
$ws.Range('D10').Value = '18.74'
$ws.Range('E10').Value = '  -5.57%  '
$ws.Range('D11').Value = '0.7791'
$ws.Range('E11').Value = '  -1.61%  '
$ws.Range('D12').Value = '0.07761'
$ws.Range('E12').Value = '  +0.48%  '
$ws.Range('D13').Value = '1.841.76'
$ws.Range('E13').Value = '  -1.80%  '
$ws.Range('D14').Value = '87.98'
$ws.Range('E14').Value = '  -1.84%  '
$ws.Range('D15').Value = '5.011'
$ws.Range('E15').Value = '  -2.24%  '
$ws.Range('D16').Value = '0.9985'
$ws.Range('E16').Value = '  -0.12%  '
$ws.Range('D17').Value = '13.92'
$ws.Range('E17').Value = '  -3.33%  '
$ws.Range('D18').Value = '0.9999'
$ws.Range('E18').Value = '  +0.00%  '
$ws.Range('D19').Value = '0.000007970'
$ws.Range('E19').Value = '  -0.13%  '
$ws.Range('D20').Value = '26.485.66'
$ws.Range('E20').Value = '  -0.54%  '
$ws.Range('D21').Value = '2.069.07'
$ws.Range('E21').Value = '  -2.41%  '
$ws.Range('D22').Value = '4.617'
$ws.Range('E22').Value = '  -2.15%  '
$ws.Range('D23').Value = '9.550'
$ws.Range('E23').Value = '  -4.23%  '
$ws.Range('D24').Value = '5.982'
$ws.Range('E24').Value = '  -1.88%  '
$ws.Range('D25').Value = '144.66'
$ws.Range('E25').Value = '  -0.69%  '
$ws.Range('D26').Value = '2.181'
$ws.Range('E26').Value = '  -7.16%  '
$ws.Range('D27').Value = '1.649'
$ws.Range('E27').Value = '  -0.33%  '
$ws.Range('D28').Value = '17.00'
$ws.Range('E28').Value = '  -0.98%  '
$ws.Range('D29').Value = '111.81'
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('D30').Value = '4.161'
$ws.Range('E30').Value = '  -3.39%  '
$ws.Range('D31').Value = '4.128'
$ws.Range('E31').Value = '  -3.85%  '
$ws.Range('D32').Value = '0.08708'
$ws.Range('E32').Value = '  -2.07%  '
$ws.Range('D33').Value = '0.04835'
$ws.Range('E33').Value = '  -1.30%  '
$ws.Range('D34').Value = '0.7241'
$ws.Range('E34').Value = '  -0.28%  '
$ws.Range('D35').Value = '1.130'
$ws.Range('E35').Value = '  -2.32%  '
$ws.Range('E36').Value = '  -1.28%  '
$ws.Range('E37').Value = '  -3.89%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').Value = '2.237'
$ws.Range('E38').Value = '  -3.47%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '0.01780'
$ws.Range('E39').Value = '  -3.29%  '
$ws.Range('D40').Value = '0.4859'
$ws.Range('E40').Value = '  -4.16%  '
$ws.Range('D41').Value = '0.9107'
$ws.Range('E41').Value = '  -2.92%  '
$ws.Range('D42').Value = '111.28'
$ws.Range('E42').Value = '  -4.06%  '
$ws.Range('D43').Value = '6.065'
$ws.Range('E43').Value = '  -0.95%  '
$ws.Range('D44').Value = '0.9993'
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('D45').Value = '7.758'
$ws.Range('E45').Value = '  -3.01%  '
$ws.Range('D46').Value = '0.4187'
$ws.Range('E46').Value = '  -4.74%  '
$ws.Range('D47').Value = '0.05932'
$ws.Range('E47').Value = '  +0.02%  '
$ws.Range('D48').Value = '9.047'
$ws.Range('E48').Value = '  -2.59%  '
$ws.Range('D49').Value = '35.06'
$ws.Range('E49').Value = '  -2.64%  '
$ws.Range('D50').Value = '0.1231'
$ws.Range('E50').Value = '  -6.85%  '
$ws.Range('D51').Value = '0.8879'
$ws.Range('E51').Value = '  +1.36%  '

$ws.Range('D2:D51').Style = 'Normal'
